$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.374.52'
$ws.Range('E2').Value = '  -1.26%  '
$ws.Range('D3').Value = '1.891.67'
$ws.Range('E3').Value = '  -1.41%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '238.00'
$ws.Range('E5').Value = '  -1.46%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.0000'
$ws.Range('E6').Value = '  -0.03%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4828'
$ws.Range('E7').Value = '  -2.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2900'
$ws.Range('E8').Value = '  -4.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06605'
$ws.Range('E9').Value = '  -2.90%  '
$ws.Range('D10').Value = '1.884.33'
$ws.Range('E10').Value = '  -1.75%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '16.89'
$ws.Range('E11').Value = '  -2.37%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07407'
$ws.Range('E12').Value = '  +0.95%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.167'
$ws.Range('E13').Value = '  -0.99%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '87.56'
$ws.Range('E14').Value = '  -1.95%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6621'
$ws.Range('E15').Value = '  -2.36%  '
$ws.Range('D16').Value = '30.346.63'
$ws.Range('E16').Value = '  -1.25%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.40'
$ws.Range('E17').Value = '  -2.00%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000007760'
$ws.Range('E18').Value = '  -3.28%  '
$ws.Range('E19').Value = '  -0.06%  '
$ws.Range('D20').Value = '2.153.60'
$ws.Range('E20').Value = '  -0.47%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.403'
$ws.Range('E21').Value = '  +0.20%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.000'
$ws.Range('E22').Value = '  -0.09%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '193.15'
$ws.Range('E23').Value = '  -4.40%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.167'
$ws.Range('E24').Value = '  -2.93%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.400'
$ws.Range('E25').Value = '  -3.40%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '163.65'
$ws.Range('E26').Value = '  +1.40%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.24'
$ws.Range('E27').Value = '  -3.51%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.948'
$ws.Range('E28').Value = '  -1.33%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.447'
$ws.Range('E29').Value = '  -0.86%  '
$ws.Range('E30').Value = '  -2.37%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09131'
$ws.Range('E31').Value = '  -0.78%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.049'
$ws.Range('E32').Value = '  -1.58%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05110'
$ws.Range('E33').Value = '  -4.20%  '
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.151'
$ws.Range('E34').Value = '  +1.78%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7320'
$ws.Range('E35').Value = '  -2.62%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.708'
$ws.Range('E36').Value = '  +0.33%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01792'
$ws.Range('E37').Value = '  -4.22%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.645'
$ws.Range('E38').Value = '  -3.09%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.9177'
$ws.Range('E39').Value = '  -1.50%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.084'
$ws.Range('E40').Value = '  -0.90%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.916'
$ws.Range('E41').Value = '  -1.21%  '
$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '106.43'
$ws.Range('E42').Value = '  -1.36%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4322'
$ws.Range('E43').Value = '  -4.44%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.002'
$ws.Range('E44').Value = '  -0.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '7.654'
$ws.Range('E45').Value = '  -1.84%  '
$ws.Range('E46').Value = '  -5.23%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.560'
$ws.Range('E47').Value = '  +7.22%  '
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '64.97'
$ws.Range('E48').Value = '  -11.24%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.961'
$ws.Range('E49').Value = '  -2.83%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05767'
$ws.Range('E50').Value = '  -3.33%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '33.96'
$ws.Range('E51').Value = '  -5.86%  '
